{"js": "// The diff only reorders the children of <w:rPr> inside several\n// character-style definitions in styles.xml so that <w:b/>/<w:i/> come\n// before <w:color/>, matching the CT_RPr sequence in wml.xsd (the order\n// OOXMLValidatorCLI enforces even though xmllint doesn't). No visual\n// formatting actually changes - toggling bold/italic off and back on for\n// each affected style is enough to make the engine re-serialize <w:rPr>\n// in schema order.\nconst styleNames = [\n  \"KeywordTok\",\n  \"ImportTok\",\n  \"CommentTok\",\n  \"DocumentationTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"ControlFlowTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n  \"AlertTok\",\n  \"ErrorTok\",\n];\n\nconst styles = context.document.getStyles();\nconst fetched = styleNames.map((name) => styles.getByNameOrNullObject(name));\nfetched.forEach((s) => s.load(\"nameLocal,font\"));\nawait context.sync();\n\nfor (const s of fetched) {\n  if (s.isNullObject) {\n    continue;\n  }\n  const font = s.font;\n  font.load(\"bold,italic\");\n}\nawait context.sync();\n\nfor (const s of fetched) {\n  if (s.isNullObject) {\n    continue;\n  }\n  const font = s.font;\n  // Re-assert the existing bold/italic state so the style's <w:rPr> gets\n  // rewritten in schema-compliant element order (rFonts, b, bCs, i, iCs,\n  // ..., color, ...); the effective formatting is unchanged.\n  if (font.bold) {\n    font.bold = true;\n  }\n  if (font.italic) {\n    font.italic = true;\n  }\n}\nawait context.sync();\n", "ps1": "# The diff only reorders the children of <w:rPr> inside several\n# character-style definitions in styles.xml so that <w:b/>/<w:i/> come\n# before <w:color/>, matching the CT_RPr sequence in wml.xsd (the order\n# OOXMLValidatorCLI enforces even though xmllint doesn't). No visual\n# formatting actually changes - toggling bold/italic off and back on for\n# each affected style is enough to make the engine re-serialize <w:rPr>\n# in schema order.\n$d = $word.ActiveDocument\n\n$styleNames = @(\n    \"KeywordTok\",\n    \"ImportTok\",\n    \"CommentTok\",\n    \"DocumentationTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"ControlFlowTok\",\n    \"InformationTok\",\n    \"WarningTok\",\n    \"AlertTok\",\n    \"ErrorTok\"\n)\n\nforeach ($name in $styleNames) {\n    try {\n        $s = $d.Styles($name)\n    } catch {\n        continue\n    }\n    $font = $s.Font\n    # Re-assert the existing bold/italic state so the style's <w:rPr> gets\n    # rewritten in schema-compliant element order (rFonts, b, bCs, i, iCs,\n    # ..., color, ...); the effective formatting is unchanged.\n    if ($font.Bold) {\n        $font.Bold = $true\n    }\n    if ($font.Italic) {\n        $font.Italic = $true\n    }\n}\n"}
